$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.538.13'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '1.809.11'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.581'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.25%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +6.53%  '
$ws.Range("E9").Value = '  +1.28%  '
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '2.068.08'
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.21'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("D14").Value = '1.801.13'
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.646'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.03%  '
$ws.Range("D16").Value = '34.494.94'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("D19").Value = '0.0₃0801'
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '246.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.41%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '173.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.85%  '
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("E26").Value = '  +7.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.86'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.51%  '
$ws.Range("E28").Value = '  +2.18%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0533'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.13%  '
$ws.Range("E32").Value = '  +0.61%  '
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").Value = '  +1.47%  '
$ws.Range("D36").Value = '1.396.84'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.50'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.95%  '
$ws.Range("E38").Value = '  -0.78%  '
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '83.66'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.98%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.961'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.76%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.38'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.55'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.70%  '
$ws.Range("E45").Value = '  +4.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0511'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("D48").Value = '1.967.55'
$ws.Range("E48").Value = '  +0.33%  '
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("E50").Value = '  -0.39%  '
$ws.Range("E51").Value = '  +0.10%  '
